$d = $word.ActiveDocument

$replacements = @(
    @("61÷3=20, 1", "76÷4=19, 0"),
    @("25÷2=12, 1", "38÷9=4, 2"),
    @("23÷4=5, 3", "19÷7=2, 5"),
    @("65÷4=16, 1", "40÷8=5, 0"),
    @("94÷8=11, 6", "19÷6=3, 1"),
    @("39÷7=5, 4", "53÷9=5, 8"),
    @("78÷8=9, 6", "50÷7=7, 1"),
    @("56÷9=6, 2", "81÷4=20, 1"),
    @("68÷5=13, 3", "56÷7=8, 0"),
    @("53÷6=8, 5", "41÷4=10, 1"),
    @("79÷4=19, 3", "75÷7=10, 5"),
    @("29÷7=4, 1", "88÷6=14, 4"),
    @("59÷8=7, 3", "20÷4=5, 0"),
    @("72÷7=10, 2", "80÷6=13, 2"),
    @("24÷2=12, 0", "89÷4=22, 1"),
    @("27÷4=6, 3", "89÷8=11, 1"),
    @("47÷7=6, 5", "86÷2=43, 0"),
    @("65÷7=9, 2", "14÷2=7, 0"),
    @("79÷8=9, 7", "21÷9=2, 3"),
    @("96÷4=24, 0", "11÷5=2, 1"),
    @("81÷9=9, 0", "72÷4=18, 0"),
    @("14÷3=4, 2", "83÷7=11, 6"),
    @("96÷8=12, 0", "26÷5=5, 1"),
    @("82÷9=9, 1", "38÷2=19, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
